$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOTAL GERAL")
$curFmt = "_-""R$""\ * #,##0.00_-;\-""R$""\ * #,##0.00_-;_-""R$""\ * ""-""??_-;_-@_-"
$ws.Range("C3").NumberFormat = $curFmt
$ws.Range("D3").NumberFormat = $curFmt
$ws.Range("I3:L3").NumberFormat = $curFmt
